$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 6 new rows before the existing row 22 ("Sector Distribution Details"
# heading), shifting everything from the old row 22 onward down to row 28+.
# This makes room for a new "Number of employees / Assets / Turnover" table.
$ws.Rows("22:27").Insert()

# New header row (bold, like the other section sub-headers in this sheet).
$ws.Range("B21").Value = "Number of employees"
$ws.Range("C21").Value = "Assets (local currency, unless noted otherwise)"
$ws.Range("D21").Value = "Turnover (local currency, unless noted otherwise)"
$ws.Range("B21:D21").Font.Bold = $true

# New data rows for Micro / Small / Medium / Large enterprise definitions.
$ws.Range("A22").Value = "Micro"
$ws.Range("B22").Value = "1-10"
$ws.Range("C22").Value = ""
$ws.Range("D22").Value = ""

$ws.Range("A23").Value = "Small"
$ws.Range("B23").Value = "11-50"
$ws.Range("C23").Value = ""
$ws.Range("D23").Value = ""

$ws.Range("A24").Value = "Medium"
$ws.Range("B24").Value = ">50"
$ws.Range("C24").Value = ""
$ws.Range("D24").Value = ""

$ws.Range("A25").Value = "Large"
$ws.Range("B25").Value = ""
$ws.Range("C25").Value = ""
$ws.Range("D25").Value = ""
